# Update "合肥-漫展信息" workbook: refresh "want to go" / "min ticket" counters on
# the existing rows and append the newly-scraped exhibitions to the
# "展览" (Exhibitions) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

# New rows scraped since last run.
$biyi = @{
    B = "2024-08-10"
    C = "合肥·比翼连枝国乙&代号鸢only"
    D = "长江东大街与东二环路交叉口向南300米东方摩域商业广场三楼 格律诗婚礼艺术中心(筑梦店)"
    E = "2024.08.10 09:00-08.10 22:00"
    F = 2
    G = 65
    H = "https://show.bilibili.com/platform/detail.html?id=88421"
    I = "//i1.hdslb.com/bfs/openplatform/202406/RANmYDJG1719330999721.jpeg"
}
$sss = @{
    B = "2024-08-18"
    C = "合肥·SSS第五人格only"
    D = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
    E = "2024.08.18 09:00-08.18 17:00"
    F = 5
    G = 68
    H = "https://show.bilibili.com/platform/detail.html?id=88430"
    I = "//i0.hdslb.com/bfs/openplatform/202406/a0qh8I1h1719660853555.png"
}

# Numeric refresh (column F = "想去人数", column G = "最低票价") applied to
# both sheets that list every event ("展览" and "全部类型" both carry the
# full roster of rows up to the "排球少年" entry).
$updates = @(
    @{ Row = 2;  F = 118  },
    @{ Row = 3;  F = 211  },
    @{ Row = 4;  F = 8    },
    @{ Row = 5;  F = 6577 },
    @{ Row = 7;  F = 429  },
    @{ Row = 8;  F = 132  },
    @{ Row = 9;  F = 5949; G = 61.2 },
    @{ Row = 12; F = 1237 },
    @{ Row = 13; F = 8    },
    @{ Row = 14; F = 85   },
    @{ Row = 16; F = 90   },
    @{ Row = 18; F = 343  },
    @{ Row = 19; F = 40   },
    @{ Row = 20; F = 4    }
)

function Apply-EventUpdates {
    param($ws, $f21Value)

    foreach ($u in $updates) {
        $ws.Cells.Item($u.Row, 6).Value = $u.F
        if ($u.ContainsKey("G")) {
            $ws.Cells.Item($u.Row, 7).Value = $u.G
        }
    }
    $ws.Cells.Item(21, 6).Value = $f21Value
}

function Fill-Row {
    param($ws, $rowIndex, $aValue, $data)

    $ws.Cells.Item($rowIndex, 1).Value = $aValue
    # Column A carries the bold/bordered/centered "row id" look used by every
    # other row. Cloning it from a known-formatted cell (row 2, always
    # present) keeps the style table from growing a near-duplicate xf, which
    # happens if Bold/Borders are poked individually on a brand-new cell.
    $ws.Range("A2").Copy()
    $ws.Range($ws.Cells.Item($rowIndex, 1), $ws.Cells.Item($rowIndex, 1)).PasteSpecial(-4122)

    # Column B holds a "yyyy-MM-dd" looking label that must stay plain text
    # (the source data keeps it as a string, not a real date) — force Text
    # format before the assignment so Excel's autodetect doesn't turn it
    # into a date serial, then drop the now-unneeded format again so the
    # cell matches its plain, unstyled siblings.
    $ws.Cells.Item($rowIndex, 2).NumberFormat = "@"
    $ws.Cells.Item($rowIndex, 2).Value = $data.B
    $ws.Cells.Item($rowIndex, 2).ClearFormats()

    $ws.Cells.Item($rowIndex, 3).Value = $data.C
    $ws.Cells.Item($rowIndex, 4).Value = $data.D
    $ws.Cells.Item($rowIndex, 5).Value = $data.E
    $ws.Cells.Item($rowIndex, 6).Value = $data.F
    $ws.Cells.Item($rowIndex, 7).Value = $data.G
    $ws.Cells.Item($rowIndex, 8).Value = $data.H
    $ws.Cells.Item($rowIndex, 9).Value = $data.I
}
# NOTE: this runtime's PowerShell subset does not bind named (-Foo bar)
# arguments on user-defined functions, so every call below passes
# parameters positionally.

# ---------------------------------------------------------------------------
# Sheet "展览" — currently rows 1..23 (header + 22 events), row 23 being the
# "银魂主题派对" entry. Insert "比翼连枝" above it (pushing "银魂" to row 24)
# and append "SSS第五人格" as the new last row (25).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

Apply-EventUpdates $ws1 4270
$ws1.Cells.Item(22, 6).Value = 37

$ws1.Rows.Item(23).Insert()
Fill-Row $ws1 23 22 $biyi

# The row that used to be 23 ("银魂") is now row 24; bump its counters/index.
$ws1.Cells.Item(24, 1).Value = 23
$ws1.Cells.Item(24, 6).Value = 185

Fill-Row $ws1 25 24 $sss

# ---------------------------------------------------------------------------
# Sheet "全部类型" — currently rows 1..24 (header + 23 events): row 23 is
# already "排球少年", row 24 is "银魂". Insert "比翼连枝" before "银魂"
# (pushing it to row 25) and append "SSS第五人格" as row 26.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

Apply-EventUpdates $ws4 4271
$ws4.Cells.Item(23, 6).Value = 37

$ws4.Rows.Item(24).Insert()
Fill-Row $ws4 24 23 $biyi

# The row that used to be 24 ("银魂") is now row 25; bump its counters/index.
$ws4.Cells.Item(25, 1).Value = 24
$ws4.Cells.Item(25, 6).Value = 185

Fill-Row $ws4 26 25 $sss
